$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra student rows (rows 12-23 in the original layout),
# leaving only the header + first 8 students.
$ws.Range("A12:C23").EntireRow.Delete() | Out-Null

# Shift the remaining block (header + 8 students, now rows 3-11) up by one
# row so the header lands on row 2 and data on rows 3-10.
$ws.Range("A1").EntireRow.Delete() | Out-Null

# Move the selection to match the saved state.
$ws.Range("M10").Select() | Out-Null
